$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 165 ---
$ws.Range("A165").Value = 45471.2916666667
$ws.Range("B165").Value = 0
$ws.Range("C165").Value = 2.96000003814697
$ws.Range("D165").Value = 2.96000003814697
$ws.Range("E165").Value = 2.96000003814697
$ws.Range("F165").Value = 2.96000003814697
$ws.Range("H165").Value = "XHS.MI"

# --- Row 166 ---
$ws.Range("A166").Value = 45474.5388541667
$ws.Range("B166").Value = 1000
$ws.Range("C166").Value = 2.88000011444092
$ws.Range("D166").Value = 2.88000011444092
$ws.Range("E166").Value = 2.88000011444092
$ws.Range("F166").Value = 2.88000011444092
$ws.Range("H166").Value = "XHS.MI"

# Column G holds numbers-as-text (shared string), matching the existing sheet
# convention. Write the text in a scratch cell far away (forced to Text
# format so Excel keeps it as a string instead of coercing to a number),
# copy only its value into the target cell, then clear the scratch cell.
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "2.96000003814697"
$ws.Range("ZZ1").Copy()
$ws.Range("G165").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "2.88000011444092"
$ws.Range("ZZ1").Copy()
$ws.Range("G166").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# --- Date formatting & font to match the existing date column (style used by A2:A164) ---
$ws.Range("A164").Copy()
$ws.Range("A165").PasteSpecial(-4122)
$ws.Range("A166").PasteSpecial(-4122)
